$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D (and E for safety) to Text format so numeric-looking
# strings (e.g. '1.00', '0.0811') are preserved exactly as text,
# matching the inlineStr cells in the target workbook.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = '42.914.37'
$ws.Range('E2').Value = '  -1.31%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = '2.564.41'
$ws.Range('E3').Value = '  -0.81%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '302.05'
$ws.Range('E5').Value = '  +0.77%  '

$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = '92.34'
$ws.Range('E6').Value = '  -3.06%  '

$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value = '0.574'
$ws.Range('E7').Value = '  +0.40%  '

$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.05%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '0.546'
$ws.Range('E9').Value = '  -0.05%  '

$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').Value = '36.17'
$ws.Range('E10').Value = '  -0.94%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.0811'
$ws.Range('E11').Value = '  +0.76%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '7.67'
$ws.Range('E12').Value = '  -0.31%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.114'
$ws.Range('E13').Value = '  +7.07%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '2.557.90'
$ws.Range('E14').Value = '  -1.32%  '

$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '0.881'
$ws.Range('E15').Value = '  +0.51%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '14.20'
$ws.Range('E16').Value = '  +0.36%  '

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '42.961.23'
$ws.Range('E17').Value = '  -1.15%  '

$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.0₃0997'
$ws.Range('E18').Value = '  +3.63%  '

$ws.Range('B19').Value = 'InternetComputer(DFINITY)'
$ws.Range('C19').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D19').Value = '12.73'
$ws.Range('E19').Value = '  +4.65%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '6.55'
$ws.Range('E20').Value = '  +0.45%  '

$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').Value = '71.53'
$ws.Range('E21').Value = '  -1.51%  '

$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = '253.23'
$ws.Range('E22').Value = '  -3.39%  '

$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').Value = '2.94'
$ws.Range('E23').Value = '  +1.27%  '

$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = '2.12'
$ws.Range('E24').Value = '  -2.28%  '

$ws.Range('B25').Value = 'EthereumClassic'
$ws.Range('C25').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D25').Value = '28.66'
$ws.Range('E25').Value = '  -0.81%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.11%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '10.25'
$ws.Range('E27').Value = '  +1.71%  '

$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').Value = '37.09'
$ws.Range('E28').Value = '  -0.04%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.13'
$ws.Range('E29').Value = '  -3.74%  '

$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '6.06'
$ws.Range('E30').Value = '  +3.07%  '

$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '153.51'
$ws.Range('E31').Value = '  +1.70%  '

$ws.Range('B32').Value = 'WEMIXToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').Value = '2.75'
$ws.Range('E32').Value = '  -1.47%  '

$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '2.15'
$ws.Range('E33').Value = '  -2.06%  '

$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '3.37'
$ws.Range('E34').Value = '  -4.91%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.0801'
$ws.Range('E35').Value = '  +0.51%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '0.114'
$ws.Range('E36').Value = '  -1.92%  '

$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').Value = '17.89'
$ws.Range('E37').Value = '  +9.90%  '

$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').Value = '0.120'
$ws.Range('E38').Value = '  +1.16%  '

$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = '23.15'
$ws.Range('E39').Value = '  -4.75%  '

$ws.Range('B40').Value = 'ApeXProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D40').Value = '2.07'
$ws.Range('E40').Value = '  +29.81%  '

$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D41').Value = '3.41'
$ws.Range('E41').Value = '  +0.02%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '0.0310'
$ws.Range('E42').Value = '  -0.11%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '3.87'
$ws.Range('E43').Value = '  +2.27%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.100.02'
$ws.Range('E44').Value = '  +1.63%  '

$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  +0.21%  '

$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '9.25'
$ws.Range('E46').Value = '  +1.74%  '

$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').Value = '85.21'
$ws.Range('E47').Value = '  -2.33%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '106.56'
$ws.Range('E48').Value = '  +1.89%  '

$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.811.76'
$ws.Range('E49').Value = '  -1.61%  '

$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').Value = '75.15'
$ws.Range('E50').Value = '  +10.43%  '

$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '1.68'
$ws.Range('E51').Value = '  +3.58%  '
